$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.135.22'
$ws.Range('E2').Value = '  +5.55%  '
$ws.Range('D3').Value = '3.506.75'
$ws.Range('E3').Value = '  +3.06%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '419.62'
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.16'
$ws.Range('E6').Value = '  +2.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.648'
$ws.Range('E7').Value = '  +4.65%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.777'
$ws.Range('E9').Value = '  +7.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  +19.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '43.14'
$ws.Range('E11').Value = '  +1.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000261'
$ws.Range('E12').Value = '  +22.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.88'
$ws.Range('E13').Value = '  +8.47%  '
$ws.Range('D14').Value = '4.060.89'
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.42'
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').Value = '3.522.82'
$ws.Range('E17').Value = '  +2.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.09'
$ws.Range('E18').Value = '  +2.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.45'
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('D20').Value = '65.039.66'
$ws.Range('E20').Value = '  +5.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '462.09'
$ws.Range('E21').Value = '  -4.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '89.99'
$ws.Range('E22').Value = '  -0.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.20'
$ws.Range('E23').Value = '  -1.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.29'
$ws.Range('E24').Value = '  +2.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.35'
$ws.Range('E25').Value = '  +1.97%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.98'
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '33.96'
$ws.Range('E27').Value = '  +3.13%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.80'
$ws.Range('E28').Value = '  +5.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '12.50'
$ws.Range('E29').Value = '  +5.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.52'
$ws.Range('E30').Value = '  -1.94%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.113'
$ws.Range('E31').Value = '  +1.88%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.164'
$ws.Range('E32').Value = '  -2.19%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '39.90'
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('B34').Value = 'Dai'
$ws.Range('C34').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '57.50'
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0501'
$ws.Range('E36').Value = '  +3.62%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0714'
$ws.Range('E37').Value = '  +37.71%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.147'
$ws.Range('E38').Value = '  +10.01%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.08'
$ws.Range('E39').Value = '  +1.65%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.77'
$ws.Range('E41').Value = '  +7.49%  '
$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.49'
$ws.Range('E42').Value = '  +7.60%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '146.43'
$ws.Range('E43').Value = '  -1.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.30'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.314'
$ws.Range('E45').Value = '  -2.05%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.01'
$ws.Range('E46').Value = '  -1.76%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.32'
$ws.Range('E47').Value = '  -1.01%  '
$ws.Range('B48').Value = 'Celestia'
$ws.Range('C48').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.86'
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.142'
$ws.Range('E49').Value = '  +2.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '21.50'
$ws.Range('E50').Value = '  -2.86%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '108.30'
$ws.Range('E51').Value = '  -3.83%  '
